# Update countries & provincias Spain
# Applies the 28-Aug-2020 22:42 data refresh to the "Pais" sheet.
# A handful of countries changed rank (their updated "Casos totales" count
# overtook a neighbouring row), so both the country name and the 7 metric
# columns (B:H) are rewritten for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-CountryRow {
    param(
        [int]$Row,
        [string]$Country,
        [double[]]$Values
    )
    $ws.Cells.Item($Row, 1).Value = $Country
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, 2 + $i).Value = $Values[$i]
    }
}

# Row -> (country, [Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes])
Set-CountryRow 4   "Estados Unidos"        @(6081949, 35315, 3360850, 2535573, 0, 730, 185526)
Set-CountryRow 23  "Alemania"              @(242101, 1536, 215495, 17246, 0, 1, 9360)
Set-CountryRow 43  "Guatemala"             @(72921, 1065, 60534, 9678, 0, 24, 2709)

# Afganistan / Costa Rica swap rank (Costa Rica's new total overtakes Afganistan)
Set-CountryRow 63  "Costa Rica"            @(38485, 1193, 14664, 23414, 0, 10, 407)
Set-CountryRow 64  "Afganistan"            @(38140, 11, 29059, 7679, 0, 1, 1402)

Set-CountryRow 75  "Estado de Palestina"   @(21251, 574, 14291, 6815, 0, 4, 145)
Set-CountryRow 87  "Sudan"                 @(13082, 37, 6603, 5656, 0, 0, 823)
Set-CountryRow 95  "Guinea"                @(9251, 38, 8282, 910, 0, 1, 59)

# Tayikistan / Gabon swap rank (Gabon's new total overtakes Tayikistan)
Set-CountryRow 98  "Gabon"                 @(8505, 37, 7187, 1265, 0, 0, 53)
Set-CountryRow 99  "Tayikistan"            @(8481, 32, 7276, 1137, 0, 0, 68)

Set-CountryRow 105 "Namibia"               @(6906, 194, 2789, 4052, 0, 5, 65)
Set-CountryRow 106 "Zimbabue"              @(6388, 96, 5043, 1150, 0, 6, 195)

# Hungria / Malaui swap rank (Malaui's new total overtakes Hungria)
Set-CountryRow 107 "Malaui"                @(5523, 27, 3143, 2206, 0, 1, 174)
Set-CountryRow 108 "Hungria"               @(5511, 132, 3759, 1138, 0, 0, 614)

Set-CountryRow 110 "Guinea Ecuatorial"     @(4941, 13, 3884, 974, 0, 0, 83)
Set-CountryRow 115 "Suazilandia"           @(4461, 28, 3210, 1160, 0, 2, 91)

# Aruba jumps ahead of Jamaica and Jordania (its new total overtakes both)
Set-CountryRow 143 "Aruba"                 @(1906, 58, 743, 1154, 0, 1, 9)
Set-CountryRow 144 "Jamaica"               @(1870, 66, 846, 1005, 0, 0, 19)
Set-CountryRow 145 "Jordania"              @(1869, 68, 1367, 487, 0, 0, 15)

Set-CountryRow 160 "Principado de Andorra" @(1124, 26, 902, 169, 0, 0, 53)
Set-CountryRow 163 "Republica del Chad"    @(1008, 4, 878, 53, 0, 0, 77)

# Refresh the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 22:42"
